$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.674.17"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").Value = "'1.629.46"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'213.43"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.489"
$ws.Range("E7").Value = "  +0.90%  "

$ws.Range("E8").Value = "  +0.79%  "

$ws.Range("E9").Value = "  +0.78%  "

$ws.Range("D10").Value = "'18.97"
$ws.Range("E10").Value = "  +4.16%  "

$ws.Range("D11").Value = "'0.0832"
$ws.Range("E11").Value = "  +2.29%  "

$ws.Range("D12").Value = "'1.857.58"
$ws.Range("E12").Value = "  +1.64%  "

$ws.Range("D13").Value = "'1.633.82"
$ws.Range("E13").Value = "  +1.62%  "

$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("E15").Value = "  +1.98%  "

$ws.Range("D16").Value = "'26.659.56"
$ws.Range("E16").Value = "  +1.57%  "

$ws.Range("D17").Value = "'62.95"
$ws.Range("E17").Value = "  +2.30%  "

$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("D20").Value = "'208.40"
$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").Value = "'4.32"
$ws.Range("E21").Value = "  +0.83%  "

$ws.Range("D22").Value = "'9.38"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("E23").Value = "  +1.57%  "

$ws.Range("E24").Value = "  -1.76%  "

$ws.Range("D25").Value = "'145.76"
$ws.Range("E25").Value = "  +1.06%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("E27").Value = "  -1.06%  "

$ws.Range("D28").Value = "'15.35"
$ws.Range("E28").Value = "  +1.05%  "

$ws.Range("D29").Value = "'6.65"
$ws.Range("E29").Value = "  +1.26%  "

$ws.Range("D30").Value = "'0.0520"
$ws.Range("E30").Value = "  +6.79%  "

$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +0.55%  "

$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  +0.51%  "

$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("E34").Value = "  +1.50%  "

$ws.Range("E35").Value = "  -0.38%  "

$ws.Range("D36").Value = "'1.165.47"
$ws.Range("E36").Value = "  +1.67%  "

$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").Value = "'0.812"
$ws.Range("E38").Value = "  +1.93%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("D41").Value = "'0.501"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("D43").Value = "'0.785"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "'1.767.13"
$ws.Range("E44").Value = "  +1.57%  "

$ws.Range("E45").Value = "  +0.48%  "

$ws.Range("E46").Value = "  +1.91%  "

$ws.Range("D47").Value = "'54.52"
$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("D48").Value = "'0.0511"
$ws.Range("E48").Value = "  +1.10%  "

$ws.Range("E49").Value = "  +4.94%  "

$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("E51").Value = "  -0.01%  "
